$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# The "Periodo Mora" (column E) entries for rows 16-21 were re-sorted from
# descending (1902, 1810, 1809, 1808, 1807, 1806) to ascending
# (1806, 1807, 1808, 1809, 1810, 1902), and the matching "Valor Mora"
# (column F) values follow the period they belong to: period 1902 carries
# 28124 while periods 1806-1810 carry 31249.

$ws.Range("E16").Value = "1806"
$ws.Range("F16").Value = 31249

$ws.Range("E17").Value = "1807"
$ws.Range("F17").Value = 31249

$ws.Range("E18").Value = "1808"
$ws.Range("F18").Value = 31249

$ws.Range("E19").Value = "1809"
$ws.Range("F19").Value = 31249

$ws.Range("E20").Value = "1810"
$ws.Range("F20").Value = 31249

$ws.Range("E21").Value = "1902"
$ws.Range("F21").Value = 28124
